$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title placeholder: the text was previously split across two runs,
#     "git " and "Tutorial" (plus an untouched leading "Welcome to " run).
#     Merge the two trailing runs into a single run "git Tutorial" while
#     leaving "Welcome to " and the line break after it untouched. ---
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$fullTitle = $titleRange.Text
$mergeStart = $fullTitle.IndexOf("git ") + 1          # 1-based COM index
$mergeLen = "git Tutorial".Length
$mergedRun = $titleRange.Characters($mergeStart, $mergeLen)
$mergedRun.Text = "git Tutorial"

# --- Subtitle placeholder: it was empty (just a trailing endParaRPr).
#     Add a new run containing "Thankyou" before that trailing mark. ---
$subtitle = $s.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange
$newRun = $subtitleRange.InsertAfter("Thankyou")
$newRun.LanguageID = "en-IN"
